$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20: new IPA script (IPAIAM0059 / OPQA-4561) ---
$ws.Range("B20").Value = "OPQA-4561"
$ws.Range("C20").Value = 'Verify that STeAM account locked message sholud be displayed in an overlay as below " Your account has been locked for 30 minutes due to too many failed attempts." <Okay>'

# --- Row 21: new IPA script (IPAIAM0060 / OPQA-4563) ---
$ws.Range("B21").Value = "OPQA-4563"
$ws.Range("C21").Value = "Verify that when STeAM account email is in an unverified status,below error message should be displayed Your email address has not yet been verified."

# --- TCID column, written after Jira id / Description so shared-string order matches ---
$ws.Range("A20").Value = "IPAIAM0059"
$ws.Range("A21").Value = "IPAIAM0060"

# --- Runmode column (reuses existing "Y" shared string) ---
$ws.Range("D20").Value = "Y"
$ws.Range("D21").Value = "Y"

# --- Copy formatting from the last existing row (19) down onto the two new rows ---
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E21").PasteSpecial(-4122)

# --- Match row height (45) used by the other wrapped-text rows ---
$ws.Rows.Item(20).RowHeight = 45
$ws.Rows.Item(21).RowHeight = 45

# --- Update the active selection to E20, matching the saved view state ---
$ws.Range("E20").Select()

$excel.CutCopyMode = $false
